$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, shifting existing rows 81:212 down to 82:213
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new record's data
$ws.Cells.Item(81, 1).Value = 9
$ws.Cells.Item(81, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(81, 3).Value = "Metropolitana"
$ws.Cells.Item(81, 4).Value = 44580
$ws.Cells.Item(81, 5).Value = 13
$ws.Cells.Item(81, 6).Value = 300000001
$ws.Cells.Item(81, 7).Value = "Rabanito"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 4300
$ws.Cells.Item(81, 11).Value = 3000
$ws.Cells.Item(81, 12).Value = 3000
$ws.Cells.Item(81, 13).Value = 3000
$ws.Cells.Item(81, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(81, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(81, 16).Value = 30
$ws.Cells.Item(81, 17).Value = 100
$ws.Cells.Item(81, 18).Value = "Hortaliza"
